$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$lo = $ws1.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Insert a new row at sheet row 69 (shifts old rows 69-109 down to 70-110)
# ---------------------------------------------------------------------
$ws1.Rows.Item(69).Insert()

# ---------------------------------------------------------------------
# 2) Append two brand new rows (111 and 112) at the bottom of the table
#    New unique strings must be introduced in this order so that the
#    shared-string table receives them as: Dbets, Gallito Parlays Exoticos,
#    Hassido (used later for the new row 69).
# ---------------------------------------------------------------------
$ws1.Range("A111").Value2 = "Dbets"
$ws1.Range("D105").Copy()
$ws1.Range("B111").PasteSpecial(-4122)
$ws1.Range("C111").PasteSpecial(-4122)
$ws1.Range("D111").PasteSpecial(-4122)
$ws1.Range("F105").Copy()
$ws1.Range("F111").PasteSpecial(-4122)
$ws1.Range("G111").PasteSpecial(-4122)
$ws1.Range("B111").Value2 = ""
$ws1.Range("C111").Value2 = ""
$ws1.Range("D111").Value2 = ""

$ws1.Range("A112").Value2 = "Gallito Parlays Exoticos"
$ws1.Range("D105").Copy()
$ws1.Range("B112").PasteSpecial(-4122)
$ws1.Range("C112").PasteSpecial(-4122)
$ws1.Range("D112").PasteSpecial(-4122)
$ws1.Range("F105").Copy()
$ws1.Range("F112").PasteSpecial(-4122)
$ws1.Range("G112").PasteSpecial(-4122)
$ws1.Range("B112").Value2 = ""
$ws1.Range("C112").Value2 = ""
$ws1.Range("D112").Value2 = ""

# ---------------------------------------------------------------------
# 3) Fill in the new row 69 (Hassido / Tipsters Europeos eu / 10000)
# ---------------------------------------------------------------------
$ws1.Range("A69").Value2 = "Hassido"
$ws1.Range("B69").Value2 = "Tipsters Europeos eu"
$ws1.Range("C69").Value2 = 10000
$ws1.Range("D105").Copy()
$ws1.Range("C69").PasteSpecial(-4122)
$ws1.Range("C69").Value2 = 10000

# ---------------------------------------------------------------------
# 4) Resize the table to the new extent
# ---------------------------------------------------------------------
$lo.Resize($ws1.Range("A1:H112"))

# ---------------------------------------------------------------------
# 5) Update view / selection state
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C69").Select()
